$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting existing rows 6-25 down to 7-26.
$ws.Rows.Item(6).Insert()

# Copy the style used for the date cell (column D) from the row below (now row 7)
# so the new row 6 date cell keeps the same custom date style.
$ws.Cells.Item(7, 4).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 6 with the new record's data.
$ws.Cells.Item(6, 1).Value = 8
$ws.Cells.Item(6, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(6, 3).Value = "Coquimbo"
$ws.Cells.Item(6, 4).Value = 44462
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100108
$ws.Cells.Item(6, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(6, 9).Value = 100108007
$ws.Cells.Item(6, 10).Value = "Coco"
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 100
$ws.Cells.Item(6, 14).Value = 19500
$ws.Cells.Item(6, 15).Value = 20000
$ws.Cells.Item(6, 16).Value = 19750
$ws.Cells.Item(6, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(6, 18).Value = "Perú"
$ws.Cells.Item(6, 19).Value = 988
$ws.Cells.Item(6, 20).Value = 20
